$d = $word.ActiveDocument

# --- Highlight the "Silph Company" paragraph and the "In the Building...Porygon..." paragraph ---
$pSilph = $d.Paragraphs.Item(20)
$pSilph.Range.HighlightColorIndex = 7   # wdYellow

$pPorygon = $d.Paragraphs.Item(21)
$pPorygon.Range.HighlightColorIndex = 7 # wdYellow

# --- Move the _GoBack bookmark from the end of the "Lucious says..." paragraph ---
# --- to the end of the "In the Building...Porygon..." paragraph ---
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$newBookmarkPos = $pPorygon.Range.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($newBookmarkPos, $newBookmarkPos))
